$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must remain literal
# text (preserve trailing zeros / exact formatting) -> force Text format first.
$textCells = @('D5', 'D6', 'D9', 'D11', 'D13', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D32', 'D33', 'D35', 'D40', 'D42', 'D44', 'D49', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.448.14'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '3.765.22'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '614.85'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').Value = '178.34'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = '3.759.78'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('D11').Value = '6.66'
$ws.Range('E11').Value = '  +5.46%  '
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('D13').Value = '40.09'
$ws.Range('E13').Value = '  -2.43%  '
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').Value = '4.397.05'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '3.766.87'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '69.501.33'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('E19').Value = '  -3.54%  '
$ws.Range('D20').Value = '509.28'
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').Value = '16.31'
$ws.Range('E21').Value = '  -3.33%  '
$ws.Range('D22').Value = '9.34'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').Value = '0.727'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = '2.53'
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').Value = '86.42'
$ws.Range('E25').Value = '  -1.64%  '
$ws.Range('D26').Value = '12.81'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('D28').Value = '10.58'
$ws.Range('E28').Value = '  -3.86%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('E31').Value = '  +2.73%  '
$ws.Range('D32').Value = '8.00'
$ws.Range('E32').Value = '  +2.98%  '
$ws.Range('D33').Value = '30.65'
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('E38').Value = '  +3.22%  '
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('D40').Value = '453.20'
$ws.Range('E40').Value = '  +8.35%  '
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('D42').Value = '49.86'
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('E43').Value = '  +6.21%  '
$ws.Range('D44').Value = '44.71'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').Value = '2.954.82'
$ws.Range('E46').Value = '  -2.80%  '
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '27.24'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '138.89'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('E51').Value = '  -1.04%  '
